# feat: add 2022-Q1 data
$wb = $excel.ActiveWorkbook
$ws2021 = $wb.Worksheets.Item("2021-Q4")
$wsTotal = $wb.Worksheets.Item("总计")

# --- 1. Create the new "2022-Q1" sheet, right before "总计" ---
# Clone "2021-Q4" so the column layout/headers/styles (fund code, name, A-index
# column, borders, bold header, etc.) are already correct, then overwrite only
# the figures that differ for 2022-Q1.
$ws2021.Copy($wsTotal)
$new2022 = $wb.ActiveSheet
$new2022.Name = "2022-Q1"

$new2022.Range("D2:G3").NumberFormat = "@"

$new2022.Range("D2").Value = "13.99"
$new2022.Range("E2").Value = "93.88"
$new2022.Range("F2").Value = "1.07"
$new2022.Range("G2").Value = "0.1497"
$new2022.Range("H2").Value = 8

$new2022.Range("D3").Value = "1.17"
$new2022.Range("E3").Value = "96.94"
$new2022.Range("F3").Value = "1.11"
$new2022.Range("G3").Value = "0.0130"
$new2022.Range("H3").Value = 8

$new2022.Range("D2:G3").Style = "Normal"

# --- 2. Insert a new top data row into "总计" for 2022-Q1 ---
# Re-fetch "总计" since copying a sheet in front of it shifts its position and
# the old $wsTotal handle would otherwise keep pointing at the copied sheet.
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("A2:D2").ClearFormats()

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.16

# Renumber the index column for the rows that shifted down
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2

# Restore the index-column style ("s=2": bold/centered/bordered) on the new row
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Keep the originally-active sheet/tab selected (copying/renaming sheets above
# moves the active tab onto the new sheet otherwise).
$wb.Worksheets.Item("2020-Q4").Activate()
